$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.730.53'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '2.664.08'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.10'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.653'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  -4.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.399'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.45%  '

$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.38%  '

$ws.Range("E14").Value = '  -2.64%  '

$ws.Range("D15").Value = '3.138.38'
$ws.Range("E15").Value = '  -0.42%  '

$ws.Range("D16").Value = '65.564.38'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").Value = '2.668.11'
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("E18").Value = '  -2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.43%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.44%  '

$ws.Range("E24").Value = '  +8.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000112'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '

$ws.Range("E27").Value = '  +1.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '563.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.09%  '

$ws.Range("E29").Value = '  -0.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.95%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("E33").Value = '  +2.94%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("E35").Value = '  -0.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.421'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '154.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.51%  '

$ws.Range("E41").Value = '  -2.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.30%  '

$ws.Range("E43").Value = '  -1.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.639'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.80%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0255'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.23%  '

$ws.Range("E50").Value = '  +2.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.804'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.82%  '
